$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "read data header bug": columns A and B were computed via
# formulas referencing a second (duplicate) header block in H:M plus a
# "Video start time" anchor cell in B9. Replace those formulas with their
# already-correct literal values, since the duplicate header data is bogus.
$rng = $ws.Range("A2:B7")
$rng.Value = $rng.Value2

# Remove the duplicated/erroneous header block and its supporting data
# that lived in columns H:M (rows 1-11).
$ws.Range("H1:M11").Delete()

# Remove the now-unused "Video start time:" label/value row.
$ws.Rows("9:9").Delete()

# Restore the active selection as last left by the editor.
[void]$ws.Range("P12").Select()
